$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 of the "Rules" sheet: the rule label in B11 changes from "R40" to "1".
# A leading apostrophe forces Excel to store the numeric-looking value as text
# (so it is written as a shared string, matching the source data type) instead
# of being interpreted as the number 1.
$ws.Range("B11").Value = "'1"
